$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2720.3845
$ws.Range("I32").Value = 1397
$ws.Range("K32").Value = 1397
$ws.Range("M32").Value = -1071
$ws.Range("H40").Value = 1937.0416
$ws.Range("I40").Value = 1820.6316
$ws.Range("K40").Value = 1820.6316
$ws.Range("M40").Value = -1645.6316

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4409.85
$ws.Range("I32").Value = 4510.3687
$ws.Range("K32").Value = 4510.3687
$ws.Range("M32").Value = -4223.3687
$ws.Range("H63").Value = 9991.579
$ws.Range("I63").Value = 10857.23
$ws.Range("K63").Value = 10857.23
$ws.Range("M63").Value = -10171.23
$ws.Range("H66").Value = 9991.579
$ws.Range("I66").Value = 10857.23
$ws.Range("K66").Value = 54286.14999999999
$ws.Range("M66").Value = -50854.14999999999
$ws.Range("H74").Value = 859.38464
$ws.Range("I74").Value = 852.1818
$ws.Range("K74").Value = 852.1818
$ws.Range("M74").Value = 21.81820000000005
$ws.Range("H77").Value = 859.38464
$ws.Range("I77").Value = 852.1818
$ws.Range("K77").Value = 4260.909
$ws.Range("M77").Value = 107.0910000000003
$ws.Range("H132").Value = 2416.7646
$ws.Range("I132").Value = 2316.9375
$ws.Range("J132").Value = 4014
$ws.Range("K132").Value = 6950.8125
$ws.Range("L132").Value = 12042
$ws.Range("M132").Value = -4420.8125
$ws.Range("N132").Value = -17102
$ws.Range("H133").Value = 500261
$ws.Range("J133").Value = 500261
$ws.Range("L133").Value = 500261
$ws.Range("N133").Value = -505321
$ws.Range("H134").Value = 84997
$ws.Range("J134").Value = 84997
$ws.Range("L134").Value = 84997
$ws.Range("N134").Value = -95137
$ws.Range("H135").Value = 550000
$ws.Range("J135").Value = 550000
$ws.Range("L135").Value = 550000
$ws.Range("N135").Value = -560140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3749.5
$ws.Range("I86").Value = 3749.5
$ws.Range("K86").Value = 3749.5
$ws.Range("M86").Value = -2626.5
$ws.Range("H89").Value = 3749.5
$ws.Range("I89").Value = 3749.5
$ws.Range("K89").Value = 18747.5
$ws.Range("M89").Value = -13131.5
$ws.Range("H134").Value = 4497
$ws.Range("I134").Value = 4497
$ws.Range("K134").Value = 13491
$ws.Range("M134").Value = -10956
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 80780
$ws.Range("J140").Value = 80780
$ws.Range("L140").Value = 80780
$ws.Range("N140").Value = -91140
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1573.2354
$ws.Range("I58").Value = 1329.4166
$ws.Range("J58").Value = 2158.4
$ws.Range("K58").Value = 1329.4166
$ws.Range("L58").Value = 2158.4
$ws.Range("M58").Value = -1126.4166
$ws.Range("N58").Value = -2564.4
$ws.Range("H99").Value = 2497.2856
$ws.Range("I99").Value = 1995.6666
$ws.Range("K99").Value = 1995.6666
$ws.Range("M99").Value = -497.6666
$ws.Range("H105").Value = 2882.5
$ws.Range("I105").Value = 1602.5
$ws.Range("K105").Value = 1602.5
$ws.Range("M105").Value = 144.5
$ws.Range("H107").Value = 1280.0555
$ws.Range("I107").Value = 1222.1666
$ws.Range("J107").Value = 1395.8334
$ws.Range("K107").Value = 1222.1666
$ws.Range("L107").Value = 1395.8334
$ws.Range("M107").Value = 697.8334
$ws.Range("N107").Value = -5235.8334
$ws.Range("H126").Value = 2497.2856
$ws.Range("I126").Value = 1995.6666
$ws.Range("K126").Value = 5986.9998
$ws.Range("M126").Value = -3516.9998
$ws.Range("H132").Value = 2596.7646
$ws.Range("I132").Value = 2774.6924
$ws.Range("J132").Value = 2018.5
$ws.Range("K132").Value = 8324.0772
$ws.Range("L132").Value = 6055.5
$ws.Range("M132").Value = -5794.0772
$ws.Range("N132").Value = -11115.5
$ws.Range("H136").Value = 1573.2354
$ws.Range("I136").Value = 1329.4166
$ws.Range("J136").Value = 2158.4
$ws.Range("K136").Value = 3988.2498
$ws.Range("L136").Value = 6475.200000000001
$ws.Range("M136").Value = -1438.2498
$ws.Range("N136").Value = -11575.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 46920692
$ws.Range("I11").Value = 53623612
$ws.Range("K11").Value = 160870836
$ws.Range("M11").Value = -160870696
$ws.Range("H37").Value = 65975.5
$ws.Range("J37").Value = 65975.5
$ws.Range("L37").Value = 197926.5
$ws.Range("N37").Value = -198150.5
$ws.Range("H97").Value = 946.3
$ws.Range("J97").Value = 495
$ws.Range("L97").Value = 1485
$ws.Range("N97").Value = -2477
$ws.Range("H128").Value = 627179.4
$ws.Range("I128").Value = 627179.4
$ws.Range("K128").Value = 1881538.2
$ws.Range("M128").Value = -1876558.2
$ws.Range("H133").Value = 14496.182
$ws.Range("I133").Value = 1891.8
$ws.Range("K133").Value = 5675.4
$ws.Range("M133").Value = -615.3999999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3342.2
$ws.Range("I122").Value = 3299.8462
$ws.Range("K122").Value = 9899.5386
$ws.Range("M122").Value = -7449.5386
$ws.Range("H132").Value = 3164.6
$ws.Range("I132").Value = 2021.1428
$ws.Range("K132").Value = 6063.428400000001
$ws.Range("M132").Value = -3533.428400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I22").Value = 2000
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -1705
$ws.Range("N22").Value = -3090
$ws.Range("I27").Value = 2000
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 2000
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -1893
$ws.Range("N27").Value = -2714
$ws.Range("H40").Value = 5448.5
$ws.Range("I40").Value = 5598
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 5598
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -5462
$ws.Range("N40").Value = -5272

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H63").Value = 49999
$ws.Range("J63").Value = 49999
$ws.Range("L63").Value = 49999
$ws.Range("N63").Value = -51247
$ws.Range("H66").Value = 49999
$ws.Range("J66").Value = 49999
$ws.Range("L66").Value = 149997
$ws.Range("N66").Value = -156237
$ws.Range("H107").Value = 812.04346
$ws.Range("I107").Value = 749.63635
$ws.Range("J107").Value = 869.25
$ws.Range("K107").Value = 2248.90905
$ws.Range("L107").Value = 2607.75
$ws.Range("M107").Value = -328.9090500000002
$ws.Range("N107").Value = -6447.75
$ws.Range("H138").Value = 84694.5
$ws.Range("J138").Value = 74999
$ws.Range("L138").Value = 74999
$ws.Range("N138").Value = -85279
